$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to be treated as text so that numeric-looking
# strings (e.g. "211.85", "1.577.58", "0.0602") are not coerced into
# floating point numbers, matching the inlineStr/text cells in the source.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.278.09"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "1.577.58"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").Value = "211.85"
$ws.Range("E5").Value = "  +2.10%  "
$ws.Range("D6").Value = "0.494"
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").Value = "22.09"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("D10").Value = "0.0602"
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("D11").Value = "0.0870"
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("D12").Value = "1.798.43"
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("D13").Value = "1.569.51"
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("E14").Value = "  +1.04%  "
$ws.Range("D15").Value = "0.522"
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("D16").Value = "27.232.24"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").Value = "62.42"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").Value = "0.0₃0704"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").Value = "216.91"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("D22").Value = "4.16"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("D23").Value = "9.28"
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("E24").Value = "  +1.19%  "
$ws.Range("D25").Value = "154.21"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("D26").Value = "6.70"
$ws.Range("E26").Value = "  +1.19%  "
$ws.Range("D27").Value = "15.13"
$ws.Range("E27").Value = "  +0.47%  "
$ws.Range("E28").Value = "  +2.53%  "
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("E30").Value = "  +3.08%  "
$ws.Range("D31").Value = "0.0474"
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("D33").Value = "3.19"
$ws.Range("E33").Value = "  +1.80%  "
$ws.Range("D34").Value = "1.457.14"
$ws.Range("E34").Value = "  +2.33%  "
$ws.Range("E35").Value = "  +5.74%  "
$ws.Range("E36").Value = "  +0.71%  "
$ws.Range("E37").Value = "  +1.26%  "
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("D39").Value = "0.537"
$ws.Range("E39").Value = "  +0.85%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "0.813"
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "5.81"
$ws.Range("E41").Value = "  +1.92%  "
$ws.Range("E42").Value = "  +0.44%  "
$ws.Range("D43").Value = "2.36"
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("D44").Value = "1.01"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").Value = "64.80"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "1.709.99"
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("D48").Value = "86.00"
$ws.Range("E48").Value = "  -1.65%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0525"
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0102"
$ws.Range("E50").Value = "  +2.16%  "
$ws.Range("D51").Value = "0.0963"
$ws.Range("E51").Value = "  +0.87%  "

# Restore the default (unstyled) cell style now that the text values are
# committed, so no stray number-format styling is left behind.
$ws.Range("D2:E51").Style = "Normal"
